$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the cells we touch so numeric-looking values (e.g. "90.386.45",
# "237.15") stay strings like the rest of the sheet, instead of Excel coercing them to numbers/dates.
foreach ($addr in @("D2","E2","D3","E3","E4","D5","E5","D6","E6","E7","D8","E8","E9","D10","E10","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","D17","D18","E18","E19","D20","E20","D21","E21","E22","D23","E23","D24","E24","D25","E25","B26","C26","D26","E26","B27","C27","D27","E27","E29","E30","E31","D32","E32","D33","E33","E34","E35","D36","E36","D37","E37","D38","E38","B39","C39","D39","E39","B40","C40","D40","E40","D41","E41","B42","C42","D42","E42","B43","C43","D43","E43","D44","E44","B46","C46","D46","E46","B47","C47","D47","E47","E48","E49","B50","C50","D50","E50","B51","C51","D51","E51")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '90.386.45'
$ws.Range('E2').Value = '  +0.31%  '
$ws.Range('D3').Value = '3.121.28'
$ws.Range('E3').Value = '  +0.94%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '237.15'
$ws.Range('E5').Value = '  +0.93%  '
$ws.Range('D6').Value = '618.26'
$ws.Range('E6').Value = '  -0.33%  '
$ws.Range('E7').Value = '  +5.26%  '
$ws.Range('D8').Value = '0.368'
$ws.Range('E8').Value = '  +2.40%  '
$ws.Range('E9').Value = '  -0.05%  '
$ws.Range('D10').Value = '0.739'
$ws.Range('E10').Value = '  +3.12%  '
$ws.Range('E11').Value = '  -22.59%  '
$ws.Range('D12').Value = '0.203'
$ws.Range('E12').Value = '  +3.15%  '
$ws.Range('D13').Value = '0.0000244'
$ws.Range('E13').Value = '  -1.20%  '
$ws.Range('D14').Value = '34.97'
$ws.Range('E14').Value = '  -0.45%  '
$ws.Range('D15').Value = '5.51'
$ws.Range('E15').Value = '  +2.40%  '
$ws.Range('D16').Value = '90.359.79'
$ws.Range('E16').Value = '  +0.52%  '
$ws.Range('D17').Value = '3.699.02'
$ws.Range('D18').Value = '3.022.68'
$ws.Range('E18').Value = '  -2.05%  '
$ws.Range('E19').Value = '  -4.78%  '
$ws.Range('D20').Value = '14.98'
$ws.Range('E20').Value = '  +8.38%  '
$ws.Range('D21').Value = '5.83'
$ws.Range('E21').Value = '  +6.89%  '
$ws.Range('E22').Value = '  -5.07%  '
$ws.Range('D23').Value = '438.95'
$ws.Range('E23').Value = '  +1.56%  '
$ws.Range('D24').Value = '9.00'
$ws.Range('E24').Value = '  +2.32%  '
$ws.Range('D25').Value = '5.93'
$ws.Range('E25').Value = '  +6.29%  '
$ws.Range('B26').Value = 'Litecoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D26').Value = '87.34'
$ws.Range('E26').Value = '  +1.44%  '
$ws.Range('B27').Value = 'Aptos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D27').Value = '11.78'
$ws.Range('E27').Value = '  -0.14%  '
$ws.Range('E29').Value = '  -0.15%  '
$ws.Range('E30').Value = '  +42.84%  '
$ws.Range('E31').Value = '  +19.00%  '
$ws.Range('D32').Value = '0.169'
$ws.Range('E32').Value = '  +8.00%  '
$ws.Range('D33').Value = '9.25'
$ws.Range('E33').Value = '  +1.82%  '
$ws.Range('E34').Value = '  +12.64%  '
$ws.Range('E35').Value = '  -7.82%  '
$ws.Range('D36').Value = '7.62'
$ws.Range('E36').Value = '  +6.98%  '
$ws.Range('D37').Value = '26.06'
$ws.Range('E37').Value = '  +1.69%  '
$ws.Range('D38').Value = '504.41'
$ws.Range('E38').Value = '  +1.27%  '
$ws.Range('B39').Value = 'Fetch.AI'
$ws.Range('C39').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D39').Value = '1.35'
$ws.Range('E39').Value = '  +6.48%  '
$ws.Range('B40').Value = 'PancakeSwap'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D40').Value = '1.92'
$ws.Range('E40').Value = '  +1.79%  '
$ws.Range('D41').Value = '0.446'
$ws.Range('E41').Value = '  +11.76%  '
$ws.Range('B42').Value = 'dogwifhat'
$ws.Range('C42').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D42').Value = '3.42'
$ws.Range('E42').Value = '  -7.90%  '
$ws.Range('B43').Value = 'MantraDAO'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D43').Value = '3.69'
$ws.Range('E43').Value = '  +0.23%  '
$ws.Range('D44').Value = '22.09'
$ws.Range('E44').Value = '  -0.04%  '
$ws.Range('B46').Value = 'ARBITRUM'
$ws.Range('C46').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D46').Value = '0.716'
$ws.Range('E46').Value = '  +5.32%  '
$ws.Range('B47').Value = 'Monero'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D47').Value = '158.47'
$ws.Range('E47').Value = '  +5.08%  '
$ws.Range('E48').Value = '  +1.88%  '
$ws.Range('E49').Value = '  +4.18%  '
$ws.Range('B50').Value = 'Filecoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D50').Value = '4.44'
$ws.Range('E50').Value = '  +1.83%  '
$ws.Range('B51').Value = 'OKB'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D51').Value = '43.87'
$ws.Range('E51').Value = '  -1.17%  '
